# Design_Doc_MyStudies.pptx edit:
# Insert a new "Updated Tentative Schedule" slide (with a Tasks/Timeframe
# table) right after the existing "Tentative Schedule" slide (position 11),
# pushing "Will I use any Datasets?", "Use Case #1/#2/#3" down by one slot.

$p = $ppt.ActivePresentation

# --- Insert the new slide at position 12, using the same "Title and
# --- Content" layout as the neighboring schedule slide.
$newSlide = $p.Slides.Add(12, 2)

# --- Title ---
$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "Updated Tentative Schedule"

# --- Replace the empty body placeholder with a 9x2 schedule table ---
$newSlide.Shapes.Item(2).Delete()

$tbl = $newSlide.Shapes.AddTable(9, 2, 66.0, 143.75, 828.0, 262.8)
$tbl.Name = "Content Placeholder 3"

$table = $tbl.Table

$rows = @(
    @("Tasks", "Timeframe", 2),
    @("Live Presentation", " Week 8", 1),
    @("Bulletin and Profile Functionality", " Week 8", 0),
    @("Flashdeck Implementation and Features", " Week 8-10", 0),
    @("NodeJs and Backend Research", "Week 8-10", 0),
    @("Presentation Slides", "Week 10", 0),
    @("Demo 2", "Week 11", 0),
    @("Backend and Database Implementation", "Weeks 11-14", 0),
    @("Final Draft Presenation Slides", "Weeks 13-15", 0)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowData = $rows[$r]
    $leftText = $rowData[0]
    $rightText = $rowData[1]
    $align = $rowData[2]

    $leftCell = $table.Cell($r + 1, 1).Shape.TextFrame.TextRange
    $leftCell.Text = $leftText
    if ($align -ne 0) {
        $leftCell.ParagraphFormat.Alignment = $align
    }

    $rightCell = $table.Cell($r + 1, 2).Shape.TextFrame.TextRange
    $rightCell.Text = $rightText
    if ($align -ne 0) {
        $rightCell.ParagraphFormat.Alignment = $align
    }
}
